$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '57.881.87'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.121.65'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '532.22'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.22'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.464'
$ws.Range('E8').Value = '  +3.38%  '
$ws.Range('B9').Value = 'Toncoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.29'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.409'
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.659.83'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.136'
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.46'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000163'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '57.928.01'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.123.68'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.99'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.58'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.99'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '350.55'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.92'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.503'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.167'
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0₃0871'
$ws.Range('E27').Value = '  -7.43%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.21'
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.87'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.05'
$ws.Range('E30').Value = '  -5.44%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '21.23'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.93'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.14'
$ws.Range('E33').Value = '  -5.43%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '158.77'
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.03'
$ws.Range('E35').Value = '  -2.38%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '25.97'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.25'
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.69'
$ws.Range('E38').Value = '  +4.28%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0669'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.01'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.696'
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.394.06'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '36.95'
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('B44').Value = 'RenzoRestakedETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.162.48'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0266'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.962'
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.03'
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '19.80'
$ws.Range('E49').Value = '  -3.00%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.738'
$ws.Range('E50').Value = '  -2.64%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0912'
$ws.Range('E51').Value = '  +1.61%  '
